# Applies the "Actualización automática 2025-09-11 08:30:08" update to the
# LOZANO MOLINA TITO workbook: a sale of 278.41 (PIEDRA SINTERIZADA, client
# "MATERIALES PARA DECORACION DECORCASA CIA. LTDA.", month "septiembre") was
# added, and the dependent totals/summary cells are refreshed accordingly.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": L18 (PIEDRA SINTERIZADA for that client) ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L18").Value = 835.22

# --- Sheet "VENTA MENSUAL": F18 (septiembre for that client) and F31 (total) ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F18").Value = 835.22
$wsMensual.Range("F31").Value = 5517.16

# --- Sheet "CUMPLIMIENTO MENSUAL": PIEDRA SINTERIZADA row (11) and TOTAL row (15) ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D11").Value = 3610.62
$wsCumplimiento.Range("E11").Value = 2233.82916370549
$wsCumplimiento.Range("F11").Value = 0.6177861931663718

$wsCumplimiento.Range("D15").Value = 5591.450000000001
$wsCumplimiento.Range("E15").Value = 26116.30990313501
$wsCumplimiento.Range("F15").Value = 0.1763432679281504
